# Presentation - added placeholder slides for objectives.
#
# Before:  ... slide5 "Technology", slide6 "Demo", slide7 "Experiences"
# After:   ... slide5 "Technology", "Objectives we met", "Future objectives",
#              "Demo" (unchanged), "Experiences" (unchanged)
#
# Two new "Title and Content" slides are inserted right before the existing
# "Demo" slide; the original "Demo"/"Experiences" slides are pushed down and
# left untouched.

$p = $ppt.ActivePresentation

# "Experiences" (position 7) already uses the same Title+Content layout as
# "Demo" and its title placeholder has a plain <a:bodyPr/> (no autofit
# element) -- exactly what the two new placeholder slides need, so duplicate
# it twice as the basis for the new slides.
$source = $p.Slides.Item(7)

$newSlide1 = $source.Duplicate()
$newSlide2 = $source.Duplicate()

# Put the two new slides right before the original "Demo" slide (position 6).
$newSlide1.MoveTo(6)
$newSlide2.MoveTo(7)

# Give the new slides their own titles; leave their (empty) content
# placeholders as-is.
$newSlide1.Shapes.Item(1).TextFrame.TextRange.Text = "Objectives we met"
$newSlide2.Shapes.Item(1).TextFrame.TextRange.Text = "Future objectives"

# The original "Demo" (now position 8) and "Experiences" (now position 9)
# slides are left completely unchanged.
